# Build site at 2023-04-12 14:53:07 UTC
# Updates the LOM3095 syllabus sheet: fills in the "Objetivos" body text,
# inserts a new row for "Docentes responsaveis" value, updates the
# "Programa resumido" / "Programa" bodies, fills "Metodo" body, updates
# "Criterio" / "Norma de recuperacao" bodies, and appends a new
# "Bibliografia" row at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (Objetivos:) was missing its B/C body text - fill it in ---
$ws.Range("B10").Value = 'Apresentar conceitos sobre fontes renováveis para geração de energia térmica, elétrica e veicular, dentre outras.'
$ws.Range("C10").Value = 'Apresentar conceitos sobre fontes renováveis para geração de energia térmica, elétrica e veicular, dentre outras.'

# --- Insert a new row 13 to hold the "Docentes responsaveis:" value
#     (the teacher name moves out of the "Objetivos:" row into its own
#     row right after the "Docentes responsaveis:" label row) ---
$ws.Rows("13:13").Insert()
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '5840521 - Rosa Ana Conte'
$ws.Range("C13").Value = '5840521 - Rosa Ana Conte'
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 14 (Programa resumido:) body text ---
$ws.Range("B14").Value = 'Fontes renováveis e tecnologias limpas para geração de energia. Estudo dos sistemas atuais nacionais e mundiais.'
$ws.Range("C14").Value = 'Fontes renováveis e tecnologias limpas para geração de energia. Estudo dos sistemas atuais nacionais e mundiais.'

# --- Row 16 (Programa:) body text ---
$ws.Range("B16").Value = 'Sistemas energéticos nacionais e mundiais: fontes renováveis e fósseis. Geração de energia por fontes renováveis: solar térmica e fotovoltaica; eólica; marítima. Geração de biomassa para fins energéticos. Gerenciamento de resíduos sólidos urbanos: recicláveis e não recicláveis; programas empresariais para logística reversa; a questão dos polímeros; reflorestamento; processamento do lixo úmido doméstico. Integração de fontes renováveis para geração de energia: ciclos térmicos híbridos'
$ws.Range("C16").Value = 'Sistemas energéticos nacionais e mundiais: fontes renováveis e fósseis. Geração de energia por fontes renováveis: solar térmica e fotovoltaica; eólica; marítima. Geração de biomassa para fins energéticos. Gerenciamento de resíduos sólidos urbanos: recicláveis e não recicláveis; programas empresariais para logística reversa; a questão dos polímeros; reflorestamento; processamento do lixo úmido doméstico. Integração de fontes renováveis para geração de energia: ciclos térmicos híbridos'

# --- Row 19 (Metodo:) body text ---
$ws.Range("B19").Value = 'Aulas expositivas, seminários.'
$ws.Range("C19").Value = 'Aulas expositivas, seminários.'

# --- Row 20 (Criterio:) body text ---
$ws.Range("B20").Value = 'Serão aplicadas duas avaliações: uma escrita (P), com peso 0,6 e um seminário com apresentação oral e entrega de material escrito com peso 0,4 (T) que comporão a nota final (NF). A nota final será calculada através da expressão: NF = (0,4 P + 0,6 T).'
$ws.Range("C20").Value = 'Serão aplicadas duas avaliações: uma escrita (P), com peso 0,6 e um seminário com apresentação oral e entrega de material escrito com peso 0,4 (T) que comporão a nota final (NF). A nota final será calculada através da expressão: NF = (0,4 P + 0,6 T).'

# --- Row 21 (Norma de recuperacao:) body text ---
$ws.Range("B21").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre e apresentada nos trabalhos, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Range("C21").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre e apresentada nos trabalhos, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'

# --- New row 22: Bibliografia: label + body text ---
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'DA Rosa, A. V. BOEKER, E., 2015, Processos de Energias Renováveis, Rio de Janeiro: Elsevier; ROCHA, J.C.; ROSA, A.H.; CARDOSO, A.A. Introdução à Química do Meio Ambiente, Porto Alegre: Bookman, 2004 GONÇALVES DA SILVA, C. De Sol a Sol, São Paulo: Oficina de Textos, 2010 ÇENGEL, Y.A.; BOLES, M.A. Thermodynamics An Engineering Approach, 6th ed., New York: McGraw Hill, 2008; BORGNAKKE, C; SONNTAG, R.E. Fundamentos da termodinâmica, São Paulo: Blucher, 2013; Vilhena, A. (coord.) Lixo Municipal: manual de gerenciamento integrado, São Paulo: CEMPRE, 3a. ed., 2010; Sítios eletrônicos sobre Reciclagem de Materiais. Cada um desses sítios tem uma apresentação de programas de reciclagem nas suas áreas de atuação e oferece ampla oferta de literatura, textos e muitas vezes, vídeos sobre reciclagem específica de certos materiais. Links úteis também são apresentados, levando a sites de empresas que apresentam seus programas de reciclagem; SUN, Y.P. Supercritical fluid technology in materials science and engineering: syntheses, properties, and applications. Marcel Dekker, 2002; MATHIAS, M.C. P.P. A FORMAÇÃO DA INDÚSTRIA GLOBAL DE GÁS NATURAL - Definição, Condicionantes e Desafios, Interciências, 2010; ABREU, F. V. BIOGÁS - Economia, regulação e sustentabilidade, Interciência, 2014; BRAND , M.A. ENERGIA DE BIOMASSA FLORESTAL, Interciência, 2010; CORRÊA, O.L.S. PETRÓLEO - Noções sobre Exploração, Perfuração, Produção e Microbiologia, Interciência, 2003.'
$ws.Range("C22").Value = 'DA Rosa, A. V. BOEKER, E., 2015, Processos de Energias Renováveis, Rio de Janeiro: Elsevier; ROCHA, J.C.; ROSA, A.H.; CARDOSO, A.A. Introdução à Química do Meio Ambiente, Porto Alegre: Bookman, 2004 GONÇALVES DA SILVA, C. De Sol a Sol, São Paulo: Oficina de Textos, 2010 ÇENGEL, Y.A.; BOLES, M.A. Thermodynamics An Engineering Approach, 6th ed., New York: McGraw Hill, 2008; BORGNAKKE, C; SONNTAG, R.E. Fundamentos da termodinâmica, São Paulo: Blucher, 2013; Vilhena, A. (coord.) Lixo Municipal: manual de gerenciamento integrado, São Paulo: CEMPRE, 3a. ed., 2010; Sítios eletrônicos sobre Reciclagem de Materiais. Cada um desses sítios tem uma apresentação de programas de reciclagem nas suas áreas de atuação e oferece ampla oferta de literatura, textos e muitas vezes, vídeos sobre reciclagem específica de certos materiais. Links úteis também são apresentados, levando a sites de empresas que apresentam seus programas de reciclagem; SUN, Y.P. Supercritical fluid technology in materials science and engineering: syntheses, properties, and applications. Marcel Dekker, 2002; MATHIAS, M.C. P.P. A FORMAÇÃO DA INDÚSTRIA GLOBAL DE GÁS NATURAL - Definição, Condicionantes e Desafios, Interciências, 2010; ABREU, F. V. BIOGÁS - Economia, regulação e sustentabilidade, Interciência, 2014; BRAND , M.A. ENERGIA DE BIOMASSA FLORESTAL, Interciência, 2010; CORRÊA, O.L.S. PETRÓLEO - Noções sobre Exploração, Perfuração, Produção e Microbiologia, Interciência, 2003.'
$ws.Rows("22:22").RowHeight = 120
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
